# code thêm tạo report lương tổng hợp
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hồ sơ nhân sự")

# Update last_edited_time column (D2:D27) to the new timestamp
$lastEditedRange = $ws.Range("D2:D27")
$lastEditedRange.Value = "2024-08-12T02:00:00.000Z"

# Row 17 (update Doanh số and Tỉ lệ đạt KPI)
$ws.Range("AO17").Value = 17000000
$ws.Range("BI17").Value = 0.5667

# Row 21 (update Doanh số and Tỉ lệ đạt KPI)
$ws.Range("AO21").Value = 12600000
$ws.Range("BI21").Value = 0.42

# Row 25 (update Doanh số and Tỉ lệ đạt KPI)
$ws.Range("AO25").Value = 20000000
$ws.Range("BI25").Value = 0.6667
